$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 353.3846
$ws.Range("I58").Value = 232.72728
$ws.Range("J58").Value = 1017
$ws.Range("K58").Value = 698.18184
$ws.Range("L58").Value = 3051
$ws.Range("M58").Value = -548.18184
$ws.Range("N58").Value = -3351

$ws.Range("H76").Value = 6186.1113
$ws.Range("I76").Value = 2900
$ws.Range("J76").Value = 7125
$ws.Range("K76").Value = 2900
$ws.Range("L76").Value = 7125
$ws.Range("M76").Value = -2585
$ws.Range("N76").Value = -7755

$ws.Range("H79").Value = 6186.1113
$ws.Range("I79").Value = 2900
$ws.Range("J79").Value = 7125
$ws.Range("K79").Value = 2900
$ws.Range("L79").Value = 7125
$ws.Range("M79").Value = -1808
$ws.Range("N79").Value = -9309

$ws.Range("H113").Value = 3237.2
$ws.Range("J113").Value = 3263.5557
$ws.Range("L113").Value = 3263.5557
$ws.Range("N113").Value = -9771.555700000001

$ws.Range("H137").Value = 4168.268
$ws.Range("I137").Value = 1410.85
$ws.Range("J137").Value = 6794.381
$ws.Range("K137").Value = 4232.549999999999
$ws.Range("L137").Value = 20383.143
$ws.Range("M137").Value = -1682.549999999999
$ws.Range("N137").Value = -25483.143

$ws.Range("H138").Value = 2110.3333
$ws.Range("I138").Value = 1351.3695
$ws.Range("J138").Value = 3107.8286
$ws.Range("K138").Value = 4054.1085
$ws.Range("L138").Value = 9323.485799999999
$ws.Range("M138").Value = 1085.8915
$ws.Range("N138").Value = -19603.4858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 6500
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = -2885
$ws.Range("N3").Value = -10230

$ws.Range("H88").Value = 3058.5881
$ws.Range("I88").Value = 2445.6667
$ws.Range("J88").Value = 3748.125
$ws.Range("K88").Value = 2445.6667
$ws.Range("L88").Value = 3748.125
$ws.Range("M88").Value = -2039.6667
$ws.Range("N88").Value = -4560.125

$ws.Range("H91").Value = 3058.5881
$ws.Range("I91").Value = 2445.6667
$ws.Range("J91").Value = 3748.125
$ws.Range("K91").Value = 2445.6667
$ws.Range("L91").Value = 3748.125
$ws.Range("M91").Value = -1041.6667
$ws.Range("N91").Value = -6556.125

$ws.Range("H102").Value = 92827.414
$ws.Range("I102").Value = 112209.11
$ws.Range("J102").Value = 34682.332
$ws.Range("K102").Value = 112209.11
$ws.Range("L102").Value = 34682.332
$ws.Range("M102").Value = -110587.11
$ws.Range("N102").Value = -37926.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 141.33333
$ws.Range("I7").Value = 149.5
$ws.Range("J7").Value = 125
$ws.Range("K7").Value = 149.5
$ws.Range("L7").Value = 125
$ws.Range("M7").Value = -36.5
$ws.Range("N7").Value = -351

$ws.Range("H86").Value = 3264.9
$ws.Range("J86").Value = 3367.875
$ws.Range("L86").Value = 3367.875
$ws.Range("N86").Value = -5613.875

$ws.Range("H89").Value = 3264.9
$ws.Range("J89").Value = 3367.875
$ws.Range("L89").Value = 16839.375
$ws.Range("N89").Value = -28071.375

$ws.Range("H107").Value = 1369.0646
$ws.Range("I107").Value = 974.8261
$ws.Range("J107").Value = 2502.5
$ws.Range("K107").Value = 974.8261
$ws.Range("L107").Value = 2502.5
$ws.Range("M107").Value = 945.1739
$ws.Range("N107").Value = -6342.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 3528.75
$ws.Range("I2").Value = 2360
$ws.Range("J2").Value = 4697.5
$ws.Range("K2").Value = 2360
$ws.Range("L2").Value = 4697.5
$ws.Range("M2").Value = -2247
$ws.Range("N2").Value = -4923.5

$ws.Range("H7").Value = 39.791668
$ws.Range("I7").Value = 34.384617
$ws.Range("J7").Value = 46.18182
$ws.Range("K7").Value = 34.384617
$ws.Range("L7").Value = 46.18182
$ws.Range("M7").Value = 78.61538300000001
$ws.Range("N7").Value = -272.18182

$ws.Range("H31").Value = 1679.61
$ws.Range("I31").Value = 994.4909
$ws.Range("J31").Value = 2516.9778
$ws.Range("K31").Value = 994.4909
$ws.Range("L31").Value = 2516.9778
$ws.Range("M31").Value = -699.4909
$ws.Range("N31").Value = -3106.9778

$ws.Range("H34").Value = 1679.61
$ws.Range("I34").Value = 994.4909
$ws.Range("J34").Value = 2516.9778
$ws.Range("K34").Value = 994.4909
$ws.Range("L34").Value = 2516.9778
$ws.Range("M34").Value = -792.4909
$ws.Range("N34").Value = -2920.9778

$ws.Range("H62").Value = 4449.967
$ws.Range("I62").Value = 4734.615
$ws.Range("J62").Value = 2599.75
$ws.Range("K62").Value = 4734.615
$ws.Range("L62").Value = 2599.75
$ws.Range("M62").Value = -4110.615
$ws.Range("N62").Value = -3847.75

$ws.Range("H65").Value = 4449.967
$ws.Range("I65").Value = 4734.615
$ws.Range("J65").Value = 2599.75
$ws.Range("K65").Value = 23673.075
$ws.Range("L65").Value = 12998.75
$ws.Range("M65").Value = -20553.075
$ws.Range("N65").Value = -19238.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 712.625
$ws.Range("I92").Value = 399.33334
$ws.Range("J92").Value = 900.6
$ws.Range("K92").Value = 1198.00002
$ws.Range("L92").Value = 2701.8
$ws.Range("M92").Value = 49.99998000000005
$ws.Range("N92").Value = -5197.8

$ws.Range("H132").Value = 4045432.8
$ws.Range("I132").Value = 12501014
$ws.Range("J132").Value = 66335.586
$ws.Range("K132").Value = 112509126
$ws.Range("L132").Value = 597020.274
$ws.Range("M132").Value = -112506596
$ws.Range("N132").Value = -602080.274

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4072.5334
$ws.Range("I122").Value = 4210
$ws.Range("J122").Value = 3797.6
$ws.Range("K122").Value = 12630
$ws.Range("L122").Value = 11392.8
$ws.Range("M122").Value = -10180
$ws.Range("N122").Value = -16292.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 482.58334
$ws.Range("I22").Value = 447.2
$ws.Range("J22").Value = 507.85715
$ws.Range("K22").Value = 447.2
$ws.Range("L22").Value = 507.85715
$ws.Range("M22").Value = -152.2
$ws.Range("N22").Value = -1097.85715

$ws.Range("H27").Value = 482.58334
$ws.Range("I27").Value = 447.2
$ws.Range("J27").Value = 507.85715
$ws.Range("K27").Value = 447.2
$ws.Range("L27").Value = 507.85715
$ws.Range("M27").Value = -340.2
$ws.Range("N27").Value = -721.85715

$ws.Range("H41").Value = 5000
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -4562
$ws.Range("N41").ClearContents()

$ws.Range("H43").Value = 20000
$ws.Range("J43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 7050.9443
$ws.Range("I136").Value = 18833.834
$ws.Range("J136").Value = 1159.5
$ws.Range("K136").Value = 56501.50199999999
$ws.Range("L136").Value = 3478.5
$ws.Range("M136").Value = -53951.50199999999
$ws.Range("N136").Value = -8578.5

